$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("A17").Value = 131063926
$ws.Range("B17").Value = 83089
$ws.Range("E17").Value = 1312
$ws.Range("F17").Value = 'Gammelgransskål'
$ws.Range("G17").Value = 'Pseudographis pinicola'
$ws.Range("H17").Value = '(Nyl.) Rehm'
$ws.Range("P17").Value = 'Torsby kommun, Vrm'
$ws.Range("Q17").Value = 408603
$ws.Range("R17").Value = 6702927
$ws.Range("S17").Value = 5
$ws.Range("AW17").Value = 'Joakim Karlsson'
$ws.Range("AX17").Value = 'Joakim Karlsson'
# Row 18
$ws.Range("A18").Value = 131066881
$ws.Range("B18").Value = 57884
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = 'Tretåig hackspett'
$ws.Range("G18").Value = 'Picoides tridactylus'
$ws.Range("H18").Value = '(Linnaeus, 1758)'
$ws.Range("K18").Value = ""
$ws.Range("L18").Value = ""
$ws.Range("M18").Value = 'äldre spår'
$ws.Range("N18").Value = ""
$ws.Range("P18").Value = 'Färntjärnen, Vrm'
$ws.Range("Q18").Value = 408720
$ws.Range("R18").Value = 6703065
$ws.Range("S18").Value = 20
$ws.Range("AC18").Value = 'Ringhack på gran'
$ws.Range("AW18").Value = 'Moa Björnberg dillner'
$ws.Range("AX18").Value = 'Moa Björnberg dillner'
# Row 19
$ws.Range("A19").Value = 131066882
$ws.Range("Q19").Value = 408723
$ws.Range("R19").Value = 6703050
# Row 20
$ws.Range("A20").Value = 131066886
$ws.Range("B20").Value = 57073
$ws.Range("D20").Value = 'LC'
$ws.Range("E20").Value = 100138
$ws.Range("F20").Value = 'Tjäder'
$ws.Range("G20").Value = 'Tetrao urogallus'
$ws.Range("H20").Value = 'Linnaeus, 1758'
$ws.Range("M20").Value = 'färska spår'
$ws.Range("Q20").Value = 408565
$ws.Range("R20").Value = 6702903
$ws.Range("AC20").ClearContents()
# Row 21
$ws.Range("A21").Value = 131066892
$ws.Range("B21").Value = 81228
$ws.Range("D21").Value = 'NT'
$ws.Range("E21").Value = 1049
$ws.Range("F21").Value = 'Kortskaftad ärgspik'
$ws.Range("G21").Value = 'Microcalicium ahlneri'
$ws.Range("H21").Value = 'Tibell'
$ws.Range("K21").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("Q21").Value = 408673
$ws.Range("R21").Value = 6702997
# Row 22
$ws.Range("A22").Value = 131066899
$ws.Range("B22").Value = 83089
$ws.Range("E22").Value = 1312
$ws.Range("F22").Value = 'Gammelgransskål'
$ws.Range("G22").Value = 'Pseudographis pinicola'
$ws.Range("H22").Value = '(Nyl.) Rehm'
$ws.Range("Q22").Value = 408724
$ws.Range("R22").Value = 6703049
# Row 23
$ws.Range("A23").Value = 131066876
$ws.Range("B23").Value = 57884
$ws.Range("E23").Value = 100109
$ws.Range("F23").Value = 'Tretåig hackspett'
$ws.Range("G23").Value = 'Picoides tridactylus'
$ws.Range("H23").Value = '(Linnaeus, 1758)'
$ws.Range("K23").Value = ""
$ws.Range("L23").Value = ""
$ws.Range("M23").Value = 'äldre spår'
$ws.Range("N23").Value = ""
$ws.Range("Q23").Value = 408772
$ws.Range("R23").Value = 6703123
$ws.Range("AC23").Value = 'Ringhack på gran'
# Row 24
$ws.Range("K24").Value = ""
$ws.Range("L24").Value = ""
$ws.Range("M24").Value = 'äldre spår'
$ws.Range("N24").Value = ""
$ws.Range("AC24").Value = 'Ringhack på gran'

Write-Output "Applied changes to rows 17-24"